{"js": "// The exercise has question/answer pairs. Two of the answers were\n// previously the placeholder \"Either\" \u2014 one following the \"How many\n// children a family will have.\" question, and one following \"The number\n// of times someone will get sick in a year.\" question. Both get replaced\n// with explanatory \"Regression is better ...\" answers. The hidden\n// \"_GoBack\" bookmark that used to sit at the end of the second \"Either\"\n// paragraph moves to the end of the (new) first paragraph.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the two target paragraphs by looking at the preceding question\n// paragraph, so this keeps working even if other \"Either\" bullets exist\n// elsewhere in the document.\nlet childrenAnswerPara = null;\nlet sickAnswerPara = null;\n\nfor (let i = 0; i < items.length - 1; i++) {\n  const current = items[i];\n  const next = items[i + 1];\n  if (current.text === \"How many children a family will have.\" && next.text === \"Either\") {\n    childrenAnswerPara = next;\n  }\n  if (current.text === \"The number of times someone will get sick in a year.\" && next.text === \"Either\") {\n    sickAnswerPara = next;\n  }\n}\n\nif (!childrenAnswerPara || !sickAnswerPara) {\n  throw new Error(\"Could not locate the expected 'Either' answer paragraphs.\");\n}\n\nconst newChildrenText =\n  \"Regression is better \\u2013 don\\u2019t want to limit the possible outcomes to, say, 5 children. \";\nconst newSickText =\n  \"Regression is better \\u2013 someone could be sick every day, and you don\\u2019t you don\\u2019t want to restrict the sample outcomes to say, 5 days. \";\n\n// The \"_GoBack\" bookmark currently lives at the end of the second\n// \"Either\" paragraph (sickAnswerPara). Remove it before rewriting text so\n// it doesn't get duplicated; it will be re-inserted on the first\n// paragraph afterwards, matching the target document.\ncontext.document.deleteBookmark(\"_GoBack\");\n\n// Replace the run text in place (keeps the existing run/paragraph\n// formatting \u2014 rFonts/color/sz/etc. \u2014 untouched).\nchildrenAnswerPara.insertText(newChildrenText, Word.InsertLocation.replace);\nsickAnswerPara.insertText(newSickText, Word.InsertLocation.replace);\nawait context.sync();\n\n// Re-insert the \"_GoBack\" bookmark at the end of the (now updated) first\n// paragraph, right after its text run.\nconst bookmarkRange = childrenAnswerPara.getRange(Word.RangeLocation.end);\nbookmarkRange.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# The exercise has question/answer pairs. Two of the answers were\n# previously the placeholder \"Either\" -- one following the \"How many\n# children a family will have.\" question, and one following \"The number\n# of times someone will get sick in a year.\" question. Both get replaced\n# with explanatory \"Regression is better ...\" answers. The hidden\n# \"_GoBack\" bookmark that used to sit at the end of the second \"Either\"\n# paragraph moves to the end of the (new) first paragraph.\n\n$d = $word.ActiveDocument\n\n$newChildrenText = \"Regression is better \" + [char]0x2013 + \" don\" + [char]0x2019 + \"t want to limit the possible outcomes to, say, 5 children. \"\n$newSickText = \"Regression is better \" + [char]0x2013 + \" someone could be sick every day, and you don\" + [char]0x2019 + \"t you don\" + [char]0x2019 + \"t want to restrict the sample outcomes to say, 5 days. \"\n\n# Locate the two target paragraphs by looking at the preceding question\n# paragraph, so this keeps working even if other \"Either\" bullets exist\n# elsewhere in the document.\n$count = $d.Paragraphs.Count\n$childrenParaIndex = -1\n$sickParaIndex = -1\nfor ($i = 1; $i -lt $count; $i++) {\n    $curText = $d.Paragraphs.Item($i).Range.Text\n    $nextText = $d.Paragraphs.Item($i + 1).Range.Text\n    if ($curText -eq \"How many children a family will have.`r\" -and $nextText -eq \"Either`r\") {\n        $childrenParaIndex = $i + 1\n    }\n    if ($curText -eq \"The number of times someone will get sick in a year.`r\" -and $nextText -eq \"Either`r\") {\n        $sickParaIndex = $i + 1\n    }\n}\n\nif ($childrenParaIndex -eq -1 -or $sickParaIndex -eq -1) {\n    throw \"Could not locate the expected 'Either' answer paragraphs.\"\n}\n\n# The \"_GoBack\" bookmark currently lives at the end of the second\n# \"Either\" paragraph. Remove it before rewriting text so it doesn't end\n# up attached to stale content; it gets re-created below on the first\n# paragraph, matching the target document.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n# Replace each answer's text (exclude the trailing paragraph mark so the\n# paragraph/run formatting -- rFonts/color/sz/etc -- stays untouched).\n$childrenRange = $d.Paragraphs.Item($childrenParaIndex).Range\n[void]$childrenRange.MoveEnd(1, -1)\n$childrenRange.Text = $newChildrenText\n\n$sickRange = $d.Paragraphs.Item($sickParaIndex).Range\n[void]$sickRange.MoveEnd(1, -1)\n$sickRange.Text = $newSickText\n\n# Re-insert the \"_GoBack\" bookmark at the end of the (now updated) first\n# paragraph's text, right after its run and before the paragraph mark.\n# A directly-collapsed Range handed to Bookmarks.Add lands at the start\n# of the document/paragraph in this host, so work around it: insert a\n# one-character placeholder at the target spot, wrap the bookmark around\n# just that character, then delete the placeholder -- the bookmark\n# collapses back down to a zero-length mark at the correct position.\n$bookmarkRange = $d.Paragraphs.Item($childrenParaIndex).Range\n[void]$bookmarkRange.MoveEnd(1, -1)\n$bookmarkRange.Collapse(0)\n$bookmarkRange.InsertAfter(\"X\")\n[void]$d.Bookmarks.Add(\"_GoBack\", $bookmarkRange)\n$goBack = $d.Bookmarks.Item(\"_GoBack\")\n$goBack.Range.Text = \"\"\n"}
